$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the previous B:I data
# (and headers) one column to the right, into C:J, matching the diff's
# "Unnamed: 0.3" header insertion ahead of the existing "Unnamed: 0.2" column.
$ws.Columns("B").Insert()

# The inserted column picks up column A's formatting for the whole column;
# the two data cells (B2:B3) should stay unstyled like the other data cells.
$ws.Range("B2:B3").ClearFormats()

# Give the new header cell B1 the same header style (bold/border/centered)
# as its neighboring header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Unnamed: 0.3"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = "Housing Purchase"
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 100000

# J2 must stay a literal text string ("1/12/2025"), not get auto-converted
# to a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "1/12/2025"
$ws.Range("J2").ClearFormats()

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "Direct Deposit"
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 0

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "1/13/2025"
$ws.Range("J3").ClearFormats()
